$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.196969696969697
$ws.Range("C2").Value = 0.5833333333333334
$ws.Range("J2").Value = 0.02272727272727273
$ws.Range("P2").Value = 0.1363636363636364
$ws.Range("S2").Value = 0.06060606060606061
$ws.Range("B3").Value = 0.01265822784810127
$ws.Range("C3").Value = 0.0189873417721519
$ws.Range("P3").Value = 0.7658227848101266
$ws.Range("S3").Value = 0.2025316455696203
$ws.Range("J4").Value = 0.06818181818181818
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.2954545454545455
$ws.Range("B6").Value = 0.075
$ws.Range("D6").Value = 0.02
$ws.Range("F6").Value = 0.045
$ws.Range("J6").Value = 0.245
$ws.Range("O6").Value = 0.015
$ws.Range("Q6").Value = 0.145
$ws.Range("R6").Value = 0.065
$ws.Range("S6").Value = 0.39
$ws.Range("B7").Value = 0.07471264367816093
$ws.Range("D7").Value = 0.01724137931034483
$ws.Range("F7").Value = 0.04597701149425287
$ws.Range("J7").Value = 0.1379310344827586
$ws.Range("Q7").Value = 0.2298850574712644
$ws.Range("R7").Value = 0.04597701149425287
$ws.Range("S7").Value = 0.4482758620689655
$ws.Range("B8").Value = 0.0660377358490566
$ws.Range("D8").Value = 0.01509433962264151
$ws.Range("F8").Value = 0.05471698113207547
$ws.Range("J8").Value = 0.1113207547169811
$ws.Range("O8").Value = 0.001886792452830189
$ws.Range("Q8").Value = 0.2113207547169811
$ws.Range("R8").Value = 0.07735849056603773
$ws.Range("S8").Value = 0.4622641509433962
$ws.Range("B9").Value = 0.06976744186046512
$ws.Range("D9").Value = 0.005813953488372093
$ws.Range("E9").Value = 0.005813953488372093
$ws.Range("F9").Value = 0.02325581395348837
$ws.Range("J9").Value = 0.1337209302325581
$ws.Range("Q9").Value = 0.186046511627907
$ws.Range("R9").Value = 0.0872093023255814
$ws.Range("S9").Value = 0.4883720930232558
$ws.Range("B10").Value = 0.09992429977289932
$ws.Range("D10").Value = 0.02271006813020439
$ws.Range("E10").Value = 0.000757002271006813
$ws.Range("F10").Value = 0.07191521574564724
$ws.Range("J10").Value = 0.129447388342165
$ws.Range("O10").Value = 0.005299015897047691
$ws.Range("Q10").Value = 0.1960635881907646
$ws.Range("R10").Value = 0.09992429977289932
$ws.Range("S10").Value = 0.3739591218773656
$ws.Range("G11").Value = 0.1623616236162362
$ws.Range("J11").Value = 0.1070110701107011
$ws.Range("K11").Value = 0.2250922509225092
$ws.Range("L11").Value = 0.4907749077490775
$ws.Range("S11").Value = 0.01476014760147601
$ws.Range("G12").Value = 0.7867647058823529
$ws.Range("J12").Value = 0.1397058823529412
$ws.Range("K12").Value = 0.007352941176470588
$ws.Range("L12").Value = 0.01470588235294118
$ws.Range("S12").Value = 0.05147058823529412
$ws.Range("G13").Value = 0.6086956521739131
$ws.Range("J13").Value = 0.2826086956521739
$ws.Range("S13").Value = 0.108695652173913
$ws.Range("F15").Value = 0.0160427807486631
$ws.Range("H15").Value = 0.1390374331550802
$ws.Range("I15").Value = 0.0748663101604278
$ws.Range("J15").Value = 0.4331550802139038
$ws.Range("K15").Value = 0.08021390374331551
$ws.Range("M15").Value = 0.0053475935828877
$ws.Range("N15").Value = 0.0053475935828877
$ws.Range("O15").Value = 0.053475935828877
$ws.Range("S15").Value = 0.1925133689839572
$ws.Range("F16").Value = 0.01142857142857143
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.06857142857142857
$ws.Range("J16").Value = 0.4571428571428571
$ws.Range("K16").Value = 0.1028571428571429
$ws.Range("M16").Value = 0.04571428571428571
$ws.Range("O16").Value = 0.03428571428571429
$ws.Range("S16").Value = 0.08
$ws.Range("F17").Value = 0.02136752136752137
$ws.Range("H17").Value = 0.2222222222222222
$ws.Range("I17").Value = 0.09401709401709402
$ws.Range("J17").Value = 0.4017094017094017
$ws.Range("K17").Value = 0.0641025641025641
$ws.Range("M17").Value = 0.02136752136752137
$ws.Range("O17").Value = 0.05982905982905983
$ws.Range("S17").Value = 0.1153846153846154
$ws.Range("F18").Value = 0.01923076923076923
$ws.Range("H18").Value = 0.1971153846153846
$ws.Range("I18").Value = 0.1105769230769231
$ws.Range("J18").Value = 0.3990384615384616
$ws.Range("K18").Value = 0.0576923076923077
$ws.Range("M18").Value = 0.03846153846153846
$ws.Range("O18").Value = 0.0673076923076923
$ws.Range("S18").Value = 0.1105769230769231
$ws.Range("F19").Value = 0.00962962962962963
$ws.Range("H19").Value = 0.2414814814814815
$ws.Range("I19").Value = 0.05925925925925926
$ws.Range("J19").Value = 0.3762962962962963
$ws.Range("K19").Value = 0.09777777777777778
$ws.Range("M19").Value = 0.01407407407407407
$ws.Range("O19").Value = 0.06518518518518518
$ws.Range("S19").Value = 0.1362962962962963
